$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the competitor data (rows 2-13) by column R (估算商户数 / estimated
# merchant count) in descending order. Ties keep their prior relative order
# (stable sort), matching a "sort by estimated merchants" report refresh.
$dataRange = $ws.Range("A2:R13")
$sortKey = $ws.Range("R2:R13")
$dataRange.Sort($sortKey, 2)

# Refresh the "data update time" column to reflect the regenerated report.
$ws.Range("P2:P13").Value = "2026-02-14 10:23:49"
